# Generate Report for Handback
#
# The 07cec6c5-5fda-43d3-be41-0e222aa8f8ed source file has now been handed
# back from localization and is in sync with en-US. Update the status on
# the Overview sheet and on each locale sheet (zh-cn, de-de), and record
# the target/handback file + handback datetime for that row.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus

# ---- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("E2").Value = "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md"
$zhcn.Range("F2").Value = "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.4989a88dc5ae437817d156855b16b3dd4195a646.zh-cn.xlf"
$zhcn.Range("G2").Value = "2016-03-10 22:44:40"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f631e8f8a59e581fe7d12d0bf84efbe66be5b9cd/e2e/07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md",
    [Type]::Missing,
    [Type]::Missing,
    "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5105eb1da0bda0db789543d64784d6262ab6f73/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/07cec6c5-5fda-43d3-be41-0e222aa8f8ed.4989a88dc5ae437817d156855b16b3dd4195a646.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.4989a88dc5ae437817d156855b16b3dd4195a646.zh-cn.xlf"
)

# ---- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = $newStatus
$dede.Range("E2").Value = "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md"
$dede.Range("F2").Value = "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.4989a88dc5ae437817d156855b16b3dd4195a646.de-de.xlf"
$dede.Range("G2").Value = "2016-03-10 22:44:57"

$dede.Hyperlinks.Add(
    $dede.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f631e8f8a59e581fe7d12d0bf84efbe66be5b9cd/e2e/07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md",
    [Type]::Missing,
    [Type]::Missing,
    "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/37b92d63750d1d871ae044cb7488a6a56a3e277b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/07cec6c5-5fda-43d3-be41-0e222aa8f8ed.4989a88dc5ae437817d156855b16b3dd4195a646.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.4989a88dc5ae437817d156855b16b3dd4195a646.de-de.xlf"
)
